$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("D3").Value = "a"
$ws.Range("E3").Value = "b"
$ws.Range("F3").Value = "c"

# Data rows D4:E15, with F4 = D4+E4 and F5:F15 a shared formula D+E
for ($i = 0; $i -lt 12; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 4).Value = $i + 1
    $ws.Cells.Item($row, 5).Value = $i + 2
}

$ws.Range("F4").Formula = "=D4+E4"
$ws.Range("F5:F15").Formula = "=D5+E5"

$ws.Range("H9").Select()
